$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update column B paths: author moved from "loren" to "nicog" Windows profile ---
$ws.Range("B1").Value = "LR"
$ws.Range("B2").Value = "C:\Users\nicog\Politecnico di Milano\DENG-SESAM - Documenti\DATASETS\Exiobase 3.8.2\MRSUT"
$ws.Range("B3").Value = "C:\Users\nicog\Politecnico di Milano\DENG-SESAM - Documenti\DATASETS\Exiobase 3.8.2\IOT"
$ws.Range("B4").Value = "C:\Users\nicog\Politecnico di Milano\DENG-SESAM - Documenti\DATASETS\Exiobase Hybrid 3.3.18"
$ws.Range("B5").Value = "C:\Users\nicog\Documents\GitHub\SESAM\GT-IOA\Database"
$ws.Range("B6").Value = "C:\Users\nicog\Documents\GitHub\SESAM\GT-IOA\Add sectors"
$ws.Range("B7").Value = "C:\Users\nicog\Documents\GitHub\SESAM\GT-IOA\Shocks"
$ws.Range("B8").Value = "C:\Users\nicog\Documents\GitHub\SESAM\GT-IOA\Results"
$ws.Range("B9").Value = "C:\Users\nicog\Documents\GitHub\SESAM\GT-IOA\Plots"
$ws.Range("B10").Value = "C:\Users\nicog\Documents\GitHub\SESAM\GT-IOA\Shocks\ShockMaster.xlsx"

# --- Add new column C ("NG") with paths pointing at the new GreenTechs project folder ---
$ws.Range("C1").Value = "NG"
$ws.Range("C2").Value = "C:\Users\nicog\Politecnico di Milano\DENG-SESAM - Documenti\DATASETS\Exiobase 3.8.2\MRSUT"
$ws.Range("C3").Value = "C:\Users\nicog\Politecnico di Milano\DENG-SESAM - Documenti\DATASETS\Exiobase 3.8.2\IOT"
$ws.Range("C4").Value = "C:\Users\nicog\Politecnico di Milano\DENG-SESAM - Documenti\DATASETS\Exiobase Hybrid 3.3.18"
$ws.Range("C5").Value = "C:\Users\nicog\Desktop\Nicolò\GitHub\GreenTechs\Database"
$ws.Range("C6").Value = "C:\Users\nicog\Desktop\Nicolò\GitHub\GreenTechs\Add sectors"
$ws.Range("C7").Value = "C:\Users\nicog\Desktop\Nicolò\GitHub\GreenTechs\Shocks"
$ws.Range("C8").Value = "C:\Users\nicog\Desktop\Nicolò\GitHub\GreenTechs\Results"
$ws.Range("C9").Value = "C:\Users\nicog\Desktop\Nicolò\GitHub\GreenTechs\Plots"
$ws.Range("C10").Value = "C:\Users\nicog\Desktop\Nicolò\GitHub\GreenTechs\ShockMaster.xlsx"

# Match the authored column width for column B
$ws.Columns.Item(2).ColumnWidth = 104.83

# Rename the cell style from the Italian-locale "Normale" to "Normal"
$wb.Styles.Add("Normal")
$wb.Styles("Normale").Delete()

# Selection mirrors the authored change (cursor left on newly-filled C2)
$ws.Range("C2").Select()
